# Add "2022-Q3" sheet with its holdings data, and add the corresponding
# summary row at the top of the "总计" (total) sheet.

$wb = $excel.ActiveWorkbook

function Set-TextCell($cell, $val) {
    # Force the cell to be stored as text even when the value looks like a
    # number (e.g. "005444", "0.84"), matching how the source data is typed.
    $cell.NumberFormat = "@"
    $cell.Value() = $val
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q3" sheet by duplicating the "2022-Q2" sheet
#    (same column layout/formatting), placed right before it.
# ---------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q2Sheet.Copy($q2Sheet)
$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# Drop the old data rows (3:7) copied from 2022-Q2, keeping only the
# header row and a single data row template (row 2).
$q3Sheet.Rows("3:7").Delete()

# Fill in the single holding row for 2022-Q3.
$q3Sheet.Cells.Item(2, 1).Value() = 0
Set-TextCell $q3Sheet.Cells.Item(2, 2) "005444"
Set-TextCell $q3Sheet.Cells.Item(2, 3) "光大保德信多策略精选18个月定期开放灵活配置混合"
Set-TextCell $q3Sheet.Cells.Item(2, 4) "0.84"
Set-TextCell $q3Sheet.Cells.Item(2, 5) "29.28"
Set-TextCell $q3Sheet.Cells.Item(2, 6) "1.54"
Set-TextCell $q3Sheet.Cells.Item(2, 7) "0.0129"
$q3Sheet.Cells.Item(2, 8).Value() = 6

# ---------------------------------------------------------------------
# 2) Update the "总计" sheet: push every existing row down by one and
#    write the new 2022-Q3 summary row at the top (row 2). Rewriting the
#    whole table directly (rather than Rows.Insert, which drags stray
#    formatting along) keeps formatting identical to the source rows.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$summaryRows = @(
    @(0, "2022-Q3", 1, 0.01),
    @(1, "2022-Q2", 6, 1.6),
    @(2, "2022-Q1", 11, 2.7),
    @(3, "2021-Q4", 3, 0.79),
    @(4, "2021-Q3", 7, 0.95),
    @(5, "2021-Q2", 10, 1.08)
)

for ($i = 0; $i -lt $summaryRows.Length; $i++) {
    $r = $i + 2
    $values = $summaryRows[$i]
    $totalSheet.Cells.Item($r, 1).Value() = $values[0]
    $totalSheet.Cells.Item($r, 2).Value() = $values[1]
    $totalSheet.Cells.Item($r, 3).Value() = $values[2]
    $totalSheet.Cells.Item($r, 4).Value() = $values[3]
}

# Row 7 is brand new territory for this sheet; give cell A7 the same
# index-column formatting (bold, centered, bordered) as the rows above it.
$totalSheet.Range("A6").Copy()
$totalSheet.Range("A7").PasteSpecial(-4122)
$totalSheet.Cells.Item(7, 1).Value() = 5

# Restore the original active sheet (sheet-copying above moves focus to
# the freshly created sheet).
$wb.Worksheets.Item("2021-Q2").Activate()
